# Loan RBI, Variable Instalments
# Insert a new (blank) column before column N on the "Repayment Schedule"
# sheet, shifting the old N/O/P columns one place to the right (-> O/P/Q),
# and update the current selection to reflect the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

$ws.Columns("N").Insert()

$ws.Range("R8").Select() | Out-Null
